# Applies the data corrections made in the "Changes into mutations method" commit
# to the distance matrix workbook ("Distancias entre las sedes").

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Fix the misspelling of "Guantánamo" (accent was on the wrong vowel) in the
# province name list.
$ws.Range("A17").Value = "Guantánamo"

# Correct the round-trip distance between Sancti Spíritus (row 9) and
# Holguín (row 15), which was recorded as 325 km but should be 400 km.
$ws.Range("O9").Value = 400
$ws.Range("I15").Value = 400

# Update the sheet's saved selection/active cell.
$ws.Range("M11").Select()
